$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("path")

# Extend the table ("Table2") by one row and fill in the new data row
$lo = $ws.ListObjects.Item("Table2")
$newRow = $lo.ListRows.Add()

$ws.Range("A10").Value = "access_db"
$ws.Range("C10").Value = "D:\OneDrive - TIQN\04.HR DB\ZK ATTENDANCE.MDB"

# Match the hyperlink-style formatting used by the other rows in column C
$ws.Range("C10").Style = $ws.Range("C9").Style

# Make "path" the active/selected sheet, with C10 selected
$ws.Select()
$ws.Range("C10").Select()
